$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell N2 (mirrors M2: thin bottom border area, no value) ---
$ws.Range("N2").Value = $null
$ws.Range("M2").Copy() | Out-Null
$ws.Range("N2").PasteSpecial(-4122) | Out-Null

# --- New year column header N3 = 2022 (mirrors M3 style) ---
$ws.Range("M3").Copy() | Out-Null
$ws.Range("N3").PasteSpecial(-4122) | Out-Null
$ws.Range("N3").Value = 2022

# --- Data rows: copy format from the corresponding M-column cell, then set value ---
$ws.Range("M4").Copy() | Out-Null
$ws.Range("N4").PasteSpecial(-4122) | Out-Null
$ws.Range("N4").Value = 9.224468514531754

$ws.Range("M5").Copy() | Out-Null
$ws.Range("N5").PasteSpecial(-4122) | Out-Null
$ws.Range("N5").Value = 4.6068543125097872

$ws.Range("M6").Copy() | Out-Null
$ws.Range("N6").PasteSpecial(-4122) | Out-Null
$ws.Range("N6").Value = 13.543910285971602

$ws.Range("M4").Copy() | Out-Null
$ws.Range("N7").PasteSpecial(-4122) | Out-Null
$ws.Range("N7").Value = 24.703327617190443

$ws.Range("M8").Copy() | Out-Null
$ws.Range("N8").PasteSpecial(-4122) | Out-Null
$ws.Range("N8").Value = 28.608474183838851

$ws.Range("M9").Copy() | Out-Null
$ws.Range("N9").PasteSpecial(-4122) | Out-Null
$ws.Range("N9").Value = 20.904451081350146

$ws.Range("M10").Copy() | Out-Null
$ws.Range("N10").PasteSpecial(-4122) | Out-Null
$ws.Range("N10").Value = 26.720095429750884

$ws.Range("M9").Copy() | Out-Null
$ws.Range("N11").PasteSpecial(-4122) | Out-Null
$ws.Range("N11").Value = 27.704327204727914

$ws.Range("M12").Copy() | Out-Null
$ws.Range("N12").PasteSpecial(-4122) | Out-Null
$ws.Range("N12").Value = 25.731792255708452

# New data columns D:N in rows 4-12, and data column N, all need "0.0" number format
$ws.Range("N4:N12").NumberFormat = "0.0"

# Update the selection to match the target end state
$ws.Range("Q5").Select() | Out-Null
